# Fruta / hortaliza, semanal
# Insert 3 new weekly price rows (Sandia, Terminal La Palmera de La Serena)
# at the top of the most-recent-first block, pushing the existing rows
# (old 110-121) down to (new 113-124).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three blank rows above row 110 (existing data shifts down).
$ws.Rows.Item(110).Insert()
$ws.Rows.Item(110).Insert()
$ws.Rows.Item(110).Insert()

# Common/fixed field values shared by every row of this market block.
$mercadoId = 8
$mercado = "Terminal La Palmera de La Serena"
$region = "Coquimbo"
$codreg = 4
$categoriaId = 100112028
$categoria = "Sandia"
$variedad = "Sin especificar"
$kgOUnidades = 1
$clasificacion = "Hortaliza"
$fecha = 44918
$unidad = "`$/unidad"
$origen = "Regi" + [char]0x00F3 + "n de O'Higgins"

# New row 110: Extra
$ws.Cells.Item(110, 1).Value = $mercadoId
$ws.Cells.Item(110, 2).Value = $mercado
$ws.Cells.Item(110, 3).Value = $region
$ws.Cells.Item(110, 4).Value = $fecha
$ws.Cells.Item(110, 5).Value = $codreg
$ws.Cells.Item(110, 6).Value = $categoriaId
$ws.Cells.Item(110, 7).Value = $categoria
$ws.Cells.Item(110, 8).Value = $variedad
$ws.Cells.Item(110, 9).Value = "Extra"
$ws.Cells.Item(110, 10).Value = 1800
$ws.Cells.Item(110, 11).Value = 3800
$ws.Cells.Item(110, 12).Value = 4000
$ws.Cells.Item(110, 13).Value = 3900
$ws.Cells.Item(110, 14).Value = $unidad
$ws.Cells.Item(110, 15).Value = $origen
$ws.Cells.Item(110, 16).Value = 3900
$ws.Cells.Item(110, 17).Value = $kgOUnidades
$ws.Cells.Item(110, 18).Value = $clasificacion

# New row 111: Primera
$ws.Cells.Item(111, 1).Value = $mercadoId
$ws.Cells.Item(111, 2).Value = $mercado
$ws.Cells.Item(111, 3).Value = $region
$ws.Cells.Item(111, 4).Value = $fecha
$ws.Cells.Item(111, 5).Value = $codreg
$ws.Cells.Item(111, 6).Value = $categoriaId
$ws.Cells.Item(111, 7).Value = $categoria
$ws.Cells.Item(111, 8).Value = $variedad
$ws.Cells.Item(111, 9).Value = "Primera"
$ws.Cells.Item(111, 10).Value = 1600
$ws.Cells.Item(111, 11).Value = 3300
$ws.Cells.Item(111, 12).Value = 3500
$ws.Cells.Item(111, 13).Value = 3400
$ws.Cells.Item(111, 14).Value = $unidad
$ws.Cells.Item(111, 15).Value = $origen
$ws.Cells.Item(111, 16).Value = 3400
$ws.Cells.Item(111, 17).Value = $kgOUnidades
$ws.Cells.Item(111, 18).Value = $clasificacion

# New row 112: Segunda
$ws.Cells.Item(112, 1).Value = $mercadoId
$ws.Cells.Item(112, 2).Value = $mercado
$ws.Cells.Item(112, 3).Value = $region
$ws.Cells.Item(112, 4).Value = $fecha
$ws.Cells.Item(112, 5).Value = $codreg
$ws.Cells.Item(112, 6).Value = $categoriaId
$ws.Cells.Item(112, 7).Value = $categoria
$ws.Cells.Item(112, 8).Value = $variedad
$ws.Cells.Item(112, 9).Value = "Segunda"
$ws.Cells.Item(112, 10).Value = 1200
$ws.Cells.Item(112, 11).Value = 2800
$ws.Cells.Item(112, 12).Value = 3000
$ws.Cells.Item(112, 13).Value = 2900
$ws.Cells.Item(112, 14).Value = $unidad
$ws.Cells.Item(112, 15).Value = $origen
$ws.Cells.Item(112, 16).Value = 2900
$ws.Cells.Item(112, 17).Value = $kgOUnidades
$ws.Cells.Item(112, 18).Value = $clasificacion
